$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BecomePartnerPage")

# Row 2: Guru99_TestBank_Header
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Guru99_TestBank_Header"
$ws.Range("C2").Value = "h2.barone"
$ws.Range("D2").Value = "CSS"

# Fill Locator Value column for rows 3-5 first
$ws.Range("C3").Value = "table tr td[align='right']"
$ws.Range("C4").Value = "input[name='emailid']"
$ws.Range("C5").Value = "input[type='submit']"

# Then fill Locator Name column for rows 5-3 (bottom-up)
$ws.Range("B5").Value = "Guru99_TestBank_Button_Submit"
$ws.Range("B4").Value = "Guru99_TestBank_TextBox_Email"
$ws.Range("B3").Value = "Guru99_TestBank_Label_Email"

# Fill remaining Sno and Locator Type for rows 3-5
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("D3").Value = "CSS"
$ws.Range("D4").Value = "CSS"
$ws.Range("D5").Value = "CSS"

# Row 6: Guru99_TestBank_AccessDetailsToDemoSite
$ws.Range("A6").Value = 5
$ws.Range("C6").Value = "//h2[text()='Access details to demo site.']"
$ws.Range("D6").Value = "Xpath"
$ws.Range("B6").Value = "Guru99_TestBank_AccessDetailsToDemoSite"

# Rows 7-11: just A column sequential numbers 6-10
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10

# Update selection to B6 as in the diff
$ws.Range("B6").Select()
